# The workbook currently has two tabs in this order: "2022-Q2", then "总计".
# Re-sort the sheet tabs so "总计" comes first, followed by "2022-Q2"
# (equivalent to dragging the "总计" tab in front of "2022-Q2" in the Excel UI).

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# Move "总计" to be immediately before "2022-Q2"
$wsTotal.Move($wsQ2)

# Keep "2022-Q2" as the active/selected tab, same as before the reorder
$wb.Worksheets.Item("2022-Q2").Activate()
